$d = $word.ActiveDocument

# Locate the paragraph that currently holds the "Method" section's opening
# text ("   " followed by the dataset-description sentence).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*The provided dataset for this task consists of four key features*") {
        $target = $p
        break
    }
}
if ($null -eq $target) {
    throw "target paragraph not found"
}

# The paragraph starts with a 3-space run ("   ") that must stay untouched.
# Replace only the remainder of the paragraph (the two sentences) in place -
# this is done as a single-paragraph XML insert so the leading run survives.
$splitPoint = $target.Range.Start + 3
$editRange = $d.Range($splitPoint, $target.Range.End)
if ($editRange.Text.Substring(0, 3) -ne "The") {
    throw "unexpected split point: [$($editRange.Text.Substring(0, 20))]"
}
$xml1 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">The provided dataset for this task consists of four key features: the date, opening price, high price of the day, and low price of the day. </w:t></w:r><w:r><w:t>This methodology emphasizes using the historical data for training, selecting the optimal model through grid search evaluated with a separate part of the training data, called evaluation set, and testing the model on unseen data without relying on the features of the test data, which are assumed to be unavailable.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$editRange.InsertXML($xml1)

# Re-resolve the (now shorter) paragraph so we know exactly where it ends.
$p1 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*without relying on the features of the test data*") {
        $p1 = $p
        break
    }
}
if ($null -eq $p1) {
    throw "rewritten paragraph not found"
}

# Insert the new "Data Preprocessing" / "Model Selection" sections right
# after that paragraph (inserted as a multi-paragraph block at a clean
# paragraph boundary, so no existing content gets swallowed).
$afterP1 = $d.Range($p1.Range.End, $p1.Range.End)
$xml2 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Kop2"/></w:pPr><w:r><w:t>Data Preprocessing</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Regardless of what model is used, the first step is to preprocess the raw data to ensure it is suitable for training and testing. Before any analysis can begin, the dataset must be inspected for any missing values or anomalies. Given that the stock market is closed on weekends and public holidays, there are weekly jumps in the data. The nature of the data is sequential, so it’s important that it is ordered chronologically, from old to new. This is critical for maintaining the time series structure and ensuring that the model only uses past data to predict future values. The task involves predicting future closing prices based on historical patterns, therefore the data was transformed into </w:t></w:r><w:r><w:t xml:space="preserve">fixed-length </w:t></w:r><w:r><w:t>sequences.</w:t></w:r><w:r><w:t xml:space="preserve"> Instead of predicting a day’s closing price, the system generates input sequences of past closing prices, which allows the model to learn patterns from historical data.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">All numerical features were normalized to a range of 0 to 1 using the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>MinMaxScaler</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>. This normalization prevents features with larger magnitude from dominating the learning process.</w:t></w:r><w:r><w:t xml:space="preserve"> Train validation split …</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Kop2"/></w:pPr><w:r><w:t>Model Selection</w:t></w:r></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>ddsdqfsfds</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p/><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$afterP1.InsertXML($xml2)
